$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "SAN DIEGO AREA TOTALS" label from B2 to A2 (keeping its formatting),
# and put a new "Totals" label in B2 with plain/default formatting.
$ws.Range("B2").Copy($ws.Range("A2"))
$ws.Range("B2").Value = "Totals"
$ws.Range("B2").Style = "Normal"

# Column A should now match column B's (bestFit) width instead of its own
# narrower custom width.
$ws.Columns.Item(1).ColumnWidth = 21.83

# Update selection to the entire column A (as in diff: sqref="A1:A1048576")
$ws.Columns.Item(1).Select()
